# Insert a new data row at row 289 (pushing the existing rows 289-308 down to
# 290-309) and populate it with the new "Femacal de La Calera" / Ají record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(289).Insert()

$ws.Range("A289").Value = 3
$ws.Range("B289").Value = "Femacal de La Calera"
$ws.Range("C289").Value = "Coquimbo"
$ws.Range("D289").Value = 44516
$ws.Range("E289").Value = 5
$ws.Range("F289").Value = 100112021
$ws.Range("G289").Value = "Ají"
$ws.Range("H289").Value = "Americana (o)"
$ws.Range("I289").Value = "Primera"
$ws.Range("J289").Value = 70
$ws.Range("K289").Value = 34000
$ws.Range("L289").Value = 35000
$ws.Range("M289").Value = 34500
$ws.Range("N289").Value = "$/caja 15 kilos"
$ws.Range("O289").Value = "Limache"
$ws.Range("P289").Value = 2300
$ws.Range("Q289").Value = 15
$ws.Range("R289").Value = "Hortaliza"
